$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; this shifts the existing rows 28-48 down to 29-49,
# matching the diff (dimension grows from A1:R48 to A1:R49).
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly data point.
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 45271
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112039
$ws.Range("G28").Value = "Ciboulette"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 150
$ws.Range("K28").Value = 2500
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 2500
$ws.Range("N28").Value = "$/docena de atados"
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 833
$ws.Range("Q28").Value = 3
$ws.Range("R28").Value = "Hortaliza"

# Column D holds dates; make sure the new cell keeps the same number format as the
# rest of the column (style index 2 in the original workbook).
$ws.Range("D28").NumberFormat = $ws.Range("D27").NumberFormat
